# ScenariosUser.xlsx - "Ajout de spécifications dans les scénarios"
#
# - Fixes a grammar typo ("des gestion" -> "de gestion") repeated on the
#   "Ajout d'une location" and "Modification d'une location" sheets.
# - Makes the "Tout les champs sont remplis" condition on the
#   "Ajout d'une location" sheet more specific (reuses the wording already
#   used on the other sheets).
# - Removes the stray "historique des locations" row on the
#   "Réservation d'une location" sheet (its table shrinks by one row).
# - Updates the selected cell on each sheet to match where the author left
#   off editing.

$wb = $excel.ActiveWorkbook

$wsDeco   = $wb.Worksheets.Item(1)   # Déconnexion
$wsAjout  = $wb.Worksheets.Item(2)   # Ajout d'une location
$wsModif  = $wb.Worksheets.Item(3)   # Modification d'une location
$wsReserv = $wb.Worksheets.Item(4)   # Réservation d'une location

# --- "Ajout d'une location" sheet -----------------------------------------
$wsAjout.Range("B7").Value2 = "Tout les champs sont remplis et les données entrées sont valides"
$wsAjout.Range("C3").Value2 = "Affiche une page de gestion des biens du user"
$wsAjout.Range("C8").Value2 = "Affiche une page de gestion des biens du user"

# --- "Modification d'une location" sheet -----------------------------------
$wsModif.Range("C3").Value2 = "Affiche une page de gestion des biens du user"
$wsModif.Range("C8").Value2 = "Affiche une page de gestion des biens du user"

# --- "Réservation d'une location" sheet ------------------------------------
# Row 6 ("Affiche une page d'historiques des locations effectuées") is
# removed entirely; everything below shifts up one row and the table /
# dimension shrink accordingly.
$wsReserv.Rows(6).Delete() | Out-Null

# --- Restore each sheet's selection -----------------------------------------
$wsAjout.Range("C20").Select() | Out-Null
$wsModif.Range("C30").Select() | Out-Null
$wsReserv.Range("C7").Select() | Out-Null
